$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.005.82"
$ws.Range("E2").Value = "  +1.09%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.896.23"
$ws.Range("E3").Value = "  +0.60%  "

# Row 4
$ws.Range("E4").Value = "  +1.69%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.13"
$ws.Range("E5").Value = "  +1.65%  "

# Row 6
$ws.Range("E6").Value = "  +1.52%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4702"
$ws.Range("E7").Value = "  -0.46%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3935"
$ws.Range("E8").Value = "  -0.90%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.75"
$ws.Range("E9").Value = "  -1.88%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08070"
$ws.Range("E10").Value = "  +0.12%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.023"
$ws.Range("E11").Value = "  -0.53%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.89"
$ws.Range("E12").Value = "  +0.20%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.901.03"
$ws.Range("E13").Value = "  +1.15%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.982"
$ws.Range("E14").Value = "  +0.13%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.144"
$ws.Range("E15").Value = "  -0.87%  "

# Row 16
$ws.Range("E16").Value = "  +1.79%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06791"
$ws.Range("E17").Value = "  +2.99%  "

# Row 18
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001057"
$ws.Range("E18").Value = "  +1.21%  "

# Row 19
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "87.68"
$ws.Range("E19").Value = "  +0.82%  "

# Row 20
$ws.Range("E20").Value = "  -0.03%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.016"
$ws.Range("E21").Value = "  +1.43%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.038.64"
$ws.Range("E22").Value = "  +1.18%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.528"
$ws.Range("E23").Value = "  +0.16%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.01"
$ws.Range("E24").Value = "  +0.08%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.349"
$ws.Range("E25").Value = "  +1.96%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.122.41"
$ws.Range("E26").Value = "  +0.93%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.71"
$ws.Range("E27").Value = "  +2.98%  "

# Row 28
$ws.Range("E28").Value = "  -0.96%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.089"
$ws.Range("E29").Value = "  -0.59%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.485"
$ws.Range("E30").Value = "  -1.85%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.13"
$ws.Range("E31").Value = "  -0.46%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9743"
$ws.Range("E32").Value = "  +1.17%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09520"
$ws.Range("E33").Value = "  -0.19%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.669"
$ws.Range("E34").Value = "  +1.31%  "

# Row 35
$ws.Range("E35").Value = "  -5.00%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.380"
$ws.Range("E36").Value = "  +1.35%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06147"
$ws.Range("E37").Value = "  +0.26%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02262"
$ws.Range("E38").Value = "  +0.15%  "

# Row 39
$ws.Range("E39").Value = "  -0.55%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.172"
$ws.Range("E40").Value = "  -0.23%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6012"
$ws.Range("E41").Value = "  -0.04%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1895"
$ws.Range("E42").Value = "  -0.28%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.37"
$ws.Range("E43").Value = "  +0.09%  "

# Row 44
$ws.Range("E44").Value = "  +1.99%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5721"
$ws.Range("E45").Value = "  +0.14%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.30"
$ws.Range("E46").Value = "  +0.73%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.410"
$ws.Range("E47").Value = "  +0.19%  "

# Row 48
$ws.Range("E48").Value = "  +0.26%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06940"
$ws.Range("E49").Value = "  +1.59%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "114.29"
$ws.Range("E50").Value = "  +3.32%  "

# Row 51
$ws.Range("E51").Value = "  +0.27%  "

